$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "divide:<type>" labels in column B (order matters for shared-string
#     table append order, matching the target diff's new string sequence) ---
$ws.Range("B9").Value2  = "divide:byte"
$ws.Range("B18").Value2 = "divide:short"
$ws.Range("B27").Value2 = "divide:int"
$ws.Range("B36").Value2 = "divide:long"
$ws.Range("B45").Value2 = "divide:float"
$ws.Range("B63").Value2 = "divide:BigInteger"
$ws.Range("B83").Value2 = "divide:IntValue"

# --- "_res_.$v$divide[...]" labels in column H ---
$ws.Range("H10").Value2 = "_res_.`$v`$divide"
$ws.Range("H19").Value2 = "_res_.`$v`$divide"
$ws.Range("H28").Value2 = "_res_.`$v`$divide"
$ws.Range("H37").Value2 = "_res_.`$v`$divide"
$ws.Range("H46").Value2 = "_res_.`$v`$divide (7)"
$ws.Range("H64").Value2 = "_res_.`$v`$divide"
$ws.Range("H84").Value2 = "_res_.`$v`$divide"

# --- Result values in column J ---
$ws.Range("J10").Value2 = 0
$ws.Range("J19").Value2 = 0
$ws.Range("J28").Value2 = 0
$ws.Range("J37").Value2 = 0
$ws.Range("J46").Value2 = 0.6666666
$ws.Range("J64").Value2 = 0
$ws.Range("J84").Value2 = 0

# --- Row 64 no longer needs the tall custom height (was ht="51") ---
$ws.Rows.Item(64).AutoFit() | Out-Null

# --- Scroll / selection state of the sheet view ---
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("O86").Select() | Out-Null
